$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Latest crypto price/volume snapshot. Price (D) values are kept as text
# (quotes stay intact even for numeric-looking prices, e.g. trailing zeros
# like "187.00" or thousand-dot formats like "26.145.85").
$updates = @(
    @{ Row = 2; D = "26.145.85"; E = "  -0.20%  " },
    @{ Row = 3; D = "1.669.12"; E = "  -0.67%  " },
    @{ Row = 4; E = "  -0.23%  " },
    @{ Row = 5; D = "210.56"; E = "  -2.64%  " },
    @{ Row = 6; D = "0.5197"; E = "  -1.30%  " },
    @{ Row = 7; D = "1.003"; E = "  -0.23%  " },
    @{ Row = 8; D = "0.2623"; E = "  -2.51%  " },
    @{ Row = 9; D = "0.06322"; E = "  -0.66%  " },
    @{ Row = 10; D = "21.16"; E = "  -1.21%  " },
    @{ Row = 11; D = "0.07540"; E = "  -1.04%  " },
    @{ Row = 12; D = "1.669.79"; E = "  -0.73%  " },
    @{ Row = 13; D = "4.442"; E = "  -1.87%  " },
    @{ Row = 14; D = "0.5487"; E = "  -4.38%  " },
    @{ Row = 15; D = "0.000008001"; E = "  -2.56%  " },
    @{ Row = 16; D = "66.32"; E = "  +0.23%  " },
    @{ Row = 17; D = "26.162.14"; E = "  -0.17%  " },
    @{ Row = 18; E = "  -0.28%  " },
    @{ Row = 19; D = "4.754"; E = "  -2.29%  " },
    @{ Row = 20; D = "187.00"; E = "  -1.54%  " },
    @{ Row = 21; D = "10.32"; E = "  -3.94%  " },
    @{ Row = 22; D = "6.210"; E = "  -0.26%  " },
    @{ Row = 23; E = "  -0.24%  " },
    @{ Row = 24; D = "150.01"; E = "  +0.57%  " },
    @{ Row = 25; D = "0.1237"; E = "  -1.76%  " },
    @{ Row = 26; D = "7.483"; E = "  -3.24%  " },
    @{ Row = 27; D = "15.81"; E = "  -0.12%  " },
    @{ Row = 28; D = "0.06300"; E = "  -0.75%  " },
    @{ Row = 29; D = "1.349"; E = "  -2.07%  " },
    @{ Row = 30; D = "1.283"; E = "  -2.44%  " },
    @{ Row = 31; D = "3.519"; E = "  -1.20%  " },
    @{ Row = 32; D = "3.409"; E = "  -4.47%  " },
    @{ Row = 33; D = "1.644"; E = "  -2.02%  " },
    @{ Row = 34; D = "1.004"; E = "  -1.68%  " },
    @{ Row = 35; D = "0.6055"; E = "  -0.78%  " },
    @{ Row = 36; E = "  -0.63%  " },
    @{ Row = 37; D = "2.763"; E = "  +0.66%  " },
    @{ Row = 38; D = "1.111.44"; E = "  +1.31%  " },
    @{ Row = 39; D = "6.108"; E = "  -1.07%  " },
    @{ Row = 40; D = "0.01613"; E = "  +0.01%  " },
    @{ Row = 41; D = "0.8652"; E = "  -2.23%  " },
    @{ Row = 42; E = "  -0.57%  " },
    @{ Row = 43; D = "100.47"; E = "  +0.12%  " },
    @{ Row = 44; D = "1.823.81"; E = "  -0.40%  " },
    @{ Row = 45; E = "  +0.51%  " },
    @{ Row = 46; D = "55.48"; E = "  -3.40%  " },
    @{ Row = 47; D = "0.9972"; E = "  -0.72%  " },
    @{ Row = 48; D = "8.061"; E = "  -0.21%  " },
    @{ Row = 49; D = "0.05238"; E = "  -0.58%  " },
    @{ Row = 50; D = "0.4243"; E = "  -0.84%  " },
    @{ Row = 51; D = "5.924"; E = "  -1.21%  " }
)

foreach ($u in $updates) {
    if ($u.ContainsKey("D")) {
        # Force text so numeric-looking strings (e.g. "187.00", "6.210")
        # keep their exact digits instead of being parsed into numbers.
        $cell = $ws.Range("D" + $u.Row)
        $cell.NumberFormat = "@"
        $cell.Value = $u.D
        $cell.ClearFormats()
    }
    if ($u.ContainsKey("E")) {
        $ws.Range("E" + $u.Row).Value = $u.E
    }
}
